$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force the cell to be treated as literal text so Excel does not
    # reinterpret numeric-looking strings (e.g. "1.00" or "26.51") as numbers.
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "75.390.40"
$ws.Range("E2").Value = "  +7.64%  "

$ws.Range("D3").Value = "2.669.09"
$ws.Range("E3").Value = "  +9.05%  "

$ws.Range("E4").Value = "  +0.02%  "

Set-TextValue "D5" "187.29"
$ws.Range("E5").Value = "  +12.08%  "

Set-TextValue "D6" "587.30"
$ws.Range("E6").Value = "  +3.37%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("E8").Value = "  +3.95%  "

Set-TextValue "D9" "0.193"
$ws.Range("E9").Value = "  +11.41%  "

$ws.Range("D10").Value = "2.670.23"
$ws.Range("E10").Value = "  +9.11%  "

$ws.Range("E11").Value = "  +1.33%  "

Set-TextValue "D12" "0.356"
$ws.Range("E12").Value = "  +6.23%  "

$ws.Range("E13").Value = "  +0.34%  "

$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").Value = "75.173.44"
$ws.Range("E14").Value = "  +7.58%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.158.39"
$ws.Range("E15").Value = "  +9.01%  "

$ws.Range("E16").Value = "  +3.75%  "

Set-TextValue "D17" "26.51"
$ws.Range("E17").Value = "  +9.77%  "

$ws.Range("D18").Value = "2.675.68"
$ws.Range("E18").Value = "  +9.29%  "

Set-TextValue "D19" "9.19"
$ws.Range("E19").Value = "  +28.50%  "

Set-TextValue "D20" "11.92"
$ws.Range("E20").Value = "  +9.85%  "

Set-TextValue "D21" "371.42"
$ws.Range("E21").Value = "  +8.85%  "

$ws.Range("E22").Value = "  +13.76%  "

Set-TextValue "D23" "4.07"
$ws.Range("E23").Value = "  +4.56%  "

Set-TextValue "D24" "6.25"
$ws.Range("E24").Value = "  +3.71%  "

Set-TextValue "D25" "1.00"
$ws.Range("E25").Value = "  +0.33%  "

$ws.Range("E26").Value = "  +5.14%  "

Set-TextValue "D27" "4.15"
$ws.Range("E27").Value = "  +8.39%  "

Set-TextValue "D28" "9.32"
$ws.Range("E28").Value = "  +9.66%  "

$ws.Range("D29").Value = "2.795.13"
$ws.Range("E29").Value = "  +8.62%  "

$ws.Range("E30").Value = "  +0.14%  "

$ws.Range("E31").Value = "  +10.24%  "

$ws.Range("E32").Value = "  +14.33%  "

Set-TextValue "D33" "520.44"
$ws.Range("E33").Value = "  +13.56%  "

Set-TextValue "D34" "7.67"
$ws.Range("E34").Value = "  +3.73%  "

$ws.Range("E35").Value = "  +7.76%  "

$ws.Range("E36").Value = "  +0.02%  "

Set-TextValue "D37" "163.32"
$ws.Range("E37").Value = "  +3.09%  "

$ws.Range("E38").Value = "  +6.03%  "

Set-TextValue "D39" "19.16"
$ws.Range("E39").Value = "  +5.06%  "

Set-TextValue "D40" "19.38"
$ws.Range("E40").Value = "  +1.47%  "

$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D42" "4.98"
$ws.Range("E42").Value = "  +12.73%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D43" "169.74"
$ws.Range("E43").Value = "  +26.00%  "

Set-TextValue "D44" "1.69"
$ws.Range("E44").Value = "  +10.48%  "

Set-TextValue "D45" "0.329"
$ws.Range("E45").Value = "  +8.45%  "

$ws.Range("E46").Value = "  +9.42%  "

$ws.Range("E47").Value = "  +11.83%  "

Set-TextValue "D48" "39.07"
$ws.Range("E48").Value = "  +2.64%  "

Set-TextValue "D49" "0.0845"
$ws.Range("E49").Value = "  +16.43%  "

$ws.Range("E50").Value = "  +7.17%  "

$ws.Range("E51").Value = "  +8.77%  "
